$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: split the run "% partecipazione (possesso)" into "% parte" /
# "cipazione (possesso)" and drop the _GoBack bookmark at the split point.
# ---------------------------------------------------------------------------
$rng1 = $d.Content
$rng1.Find.Execute("% parte", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $d.Range($rng1.End, $rng1.End)
$d.Bookmarks.Add("_GoBack", $splitPoint) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: merge the two runs " Ratio" + " diretta" into a single run
# " Ratio diretta" (same-text replace collapses the adjoining runs).
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute(" Ratio diretta", $true, $false, $false, $false, $false, $true, 1, $false, " Ratio diretta", 2)

# ---------------------------------------------------------------------------
# Change 3: the old _GoBack bookmark that used to sit after "Parametri da
# GFT per creazione Interscambio" was already relocated by the
# Bookmarks.Add call above (bookmark names are unique, so re-adding
# "_GoBack" moved the sole instance rather than creating a second one).
# Just append the two new paragraphs (plus the blank separator paragraph)
# that now follow that line.
# ---------------------------------------------------------------------------
$endRng = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRng.InsertAfter("`r`rSelezionare Partecipante: da capire`rCapire come visualizzare i saldi (se come indirette o gestite a partecipanti)")
